{"js": "// Applies the abstract/figure-legend revisions described in the commit\n// \"another round revision of abstract and figure legends\".\n//\n// Net textual changes (run-splitting in the OOXML diff does not change the\n// rendered text, so we only need to reproduce the visible wording):\n//   1. \"...growth of keystone inulin responders inferred...\"\n//      -> \"...growth of several inulin responders (e.g., Bacteroides\n//          acidifaciens nd unclassified Muribaculaceae) and their\n//          competitions inferred...\"\n//   2. \"Using a novel \"           -> \"Using a new \"\n//   3. \"exhibit positive and baseline-dependent responses to inulin\"\n//      -> \"exhibit positive and significant baseline-dependent responses to inulin\"\n//   4. \"However, the SCFA\"        -> \"Due to the baseline differences, SCFA\"\n//   5. \"were only marginally predictable...\" -> \"were, however, only marginally predictable...\"\n//   6. \" with different predictors\" -> \"\" (removed)\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacement, options) {\n  const opts = Object.assign({ matchCase: true, matchWholeWord: false }, options || {});\n  const results = body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${searchText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n\n// 1. \"keystone inulin responders ... inferred\" -> \"several inulin responders\n//    (e.g., Bacteroides acidifaciens nd unclassified Muribaculaceae) and\n//    their competitions inferred\"\nawait replaceOnce(\n  \"which can be explained by initial rapid growth of keystone inulin responders inferred from ecological network model\",\n  \"which can be explained by initial rapid growth of several inulin responders (e.g., Bacteroides acidifaciens nd unclassified Muribaculaceae) and their competitions inferred from ecological network model\"\n);\n\n// 2. \"novel\" -> \"new\"\nawait replaceOnce(\". Using a novel \", \". Using a new \");\n\n// 3. add \"significant\" before \"baseline-dependent\"\nawait replaceOnce(\n  \"exhibit positive and baseline-dependent responses to inulin\",\n  \"exhibit positive and significant baseline-dependent responses to inulin\"\n);\n\n// 4. \"However, the SCFA\" -> \"Due to the baseline differences, SCFA\"\nawait replaceOnce(\"However, the SCFA\", \"Due to the baseline differences, SCFA\");\n\n// 5. insert \", however,\" after \"were\"\nawait replaceOnce(\n  \"were only marginally predictable from microbiota composition\",\n  \"were, however, only marginally predictable from microbiota composition\"\n);\n\n// 6. remove \" with different predictors\"\nawait replaceOnce(\" with different predictors\", \"\");\n", "ps1": "# Applies the abstract/figure-legend revisions described in the commit\n# \"another round revision of abstract and figure legends\".\n#\n# Net textual changes (the OOXML diff re-splits several runs, but the\n# rendered text is what actually needs to change):\n#   1. \"...growth of keystone inulin responders inferred...\"\n#      -> \"...growth of several inulin responders (e.g., Bacteroides\n#          acidifaciens nd unclassified Muribaculaceae) and their\n#          competitions inferred...\"\n#   2. \"Using a novel \"           -> \"Using a new \"\n#   3. \"exhibit positive and baseline-dependent responses to inulin\"\n#      -> \"exhibit positive and significant baseline-dependent responses to inulin\"\n#   4. \"However, the SCFA\"        -> \"Due to the baseline differences, SCFA\"\n#   5. \"were only marginally predictable...\" -> \"were, however, only marginally predictable...\"\n#   6. \" with different predictors\" -> \"\" (removed)\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $found = $rng.Find.Execute(\n        $findText,    # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        $wdFindContinue, # Wrap\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        $wdReplaceAll # Replace\n    )\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 1. \"keystone inulin responders ... inferred\" -> \"several inulin responders\n#    (e.g., Bacteroides acidifaciens nd unclassified Muribaculaceae) and\n#    their competitions inferred\"\nReplace-Text `\n    \"which can be explained by initial rapid growth of keystone inulin responders inferred from ecological network model\" `\n    \"which can be explained by initial rapid growth of several inulin responders (e.g., Bacteroides acidifaciens nd unclassified Muribaculaceae) and their competitions inferred from ecological network model\"\n\n# 2. \"novel\" -> \"new\"\nReplace-Text \". Using a novel \" \". Using a new \"\n\n# 3. add \"significant\" before \"baseline-dependent\"\nReplace-Text `\n    \"exhibit positive and baseline-dependent responses to inulin\" `\n    \"exhibit positive and significant baseline-dependent responses to inulin\"\n\n# 4. \"However, the SCFA\" -> \"Due to the baseline differences, SCFA\"\nReplace-Text \"However, the SCFA\" \"Due to the baseline differences, SCFA\"\n\n# 5. insert \", however,\" after \"were\"\nReplace-Text `\n    \"were only marginally predictable from microbiota composition\" `\n    \"were, however, only marginally predictable from microbiota composition\"\n\n# 6. remove \" with different predictors\"\nReplace-Text \" with different predictors\" \"\"\n\nWrite-Output \"done\"\n"}
